# Update the disposition code mapping sheet so that column A (Code) uses
# the same descriptive text as column B (Description), reflecting the new
# disposition codes (the old short abbreviation codes like CONV/ACQ/DIS/...
# are no longer used).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CourtDisp")

$ws.Range("A4").Value = "Convicted"
$ws.Range("A5").Value = "Probation Without Verdict"
$ws.Range("A6").Value = "Not Guilty by Reason of Insanity"
$ws.Range("A7").Value = "Acquitted"
$ws.Range("A8").Value = "Dismissed"
$ws.Range("A9").Value = "Civil Procedure"
$ws.Range("A10").Value = "Off Calendar"
$ws.Range("A11").Value = "Guilty But Mentally Ill"
$ws.Range("A12").Value = "Transferred to Juvenile Court"
$ws.Range("A13").Value = "Mistrial"
$ws.Range("A14").Value = "Nolle Prosequi"
$ws.Range("A15").Value = "Other"
$ws.Range("A16").Value = "Extradited"
$ws.Range("A17").Value = "Not Disposition By Court"
$ws.Range("A18").Value = "Missing/Unknown"

$ws.Range("E18").Select()
